$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.255.12'
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").Value = '3.798.47'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'669.77"
$ws.Range("E5").Value = '  +6.99%  '

$ws.Range("D6").Value = "'168.50"
$ws.Range("E6").Value = '  +1.24%  '

$ws.Range("D7").Value = '3.795.84'
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("E10").Value = '  +0.31%  '

$ws.Range("D11").Value = "'0.461"
$ws.Range("E11").Value = '  +0.62%  '

$ws.Range("E12").Value = '  +4.69%  '

$ws.Range("E13").Value = '  -2.00%  '

$ws.Range("D14").Value = "'35.63"
$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").Value = '4.437.67'
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("D16").Value = '3.780.17'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '70.248.33'
$ws.Range("E17").Value = '  +1.19%  '

$ws.Range("D18").Value = "'17.65"
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("E19").Value = '  +1.19%  '

$ws.Range("D21").Value = "'11.42"
$ws.Range("E21").Value = '  +18.75%  '

$ws.Range("D22").Value = "'474.29"
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("E23").Value = '  +0.67%  '

$ws.Range("D24").Value = "'83.30"
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").Value = "'0.0000141"
$ws.Range("E25").Value = '  -4.84%  '

$ws.Range("E26").Value = '  +0.53%  '

$ws.Range("D27").Value = "'10.25"
$ws.Range("E27").Value = '  +1.99%  '

$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("D30").Value = '3.948.66'
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("E31").Value = '  +5.81%  '

$ws.Range("D32").Value = "'2.30"
$ws.Range("E32").Value = '  +1.81%  '

$ws.Range("D33").Value = "'7.39"
$ws.Range("E33").Value = '  +2.55%  '

$ws.Range("D34").Value = "'29.52"
$ws.Range("E34").Value = '  +2.27%  '

$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = '  +9.79%  '

$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'9.09"
$ws.Range("E37").Value = '  +1.12%  '

$ws.Range("D38").Value = '3.754.69'
$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("E39").Value = '  -0.31%  '

$ws.Range("E40").Value = '  -1.24%  '

$ws.Range("D41").Value = "'5.95"
$ws.Range("E41").Value = '  +2.03%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = "'0.963"
$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("E44").Value = '  +10.34%  '

$ws.Range("D46").Value = "'45.39"
$ws.Range("E46").Value = '  +5.09%  '

$ws.Range("D47").Value = "'159.14"
$ws.Range("E47").Value = '  +3.90%  '

$ws.Range("D48").Value = "'48.02"
$ws.Range("E48").Value = '  +2.71%  '

$ws.Range("D49").Value = "'0.299"
$ws.Range("E49").Value = '  +0.30%  '

$ws.Range("E50").Value = '  +4.10%  '

$ws.Range("D51").Value = "'0.000290"
$ws.Range("E51").Value = '  +4.55%  '
